# Table17.xlsx - fill in the "District of residence" (column A) labels that
# were left blank on every row except the first row of each district group,
# and grow the "Alaska / Washington / Oregon" footer group's row height to
# fit its wrapped text. Also update the sheet view's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Each group is defined by the row that already carries the district name
# (and its style) plus the list of rows below it that need the same
# district name / style repeated down column A.
$groups = @(
    @{ Header = 2;  Rows = @(3, 4, 5, 6) },
    @{ Header = 7;  Rows = @(8, 9, 10) },
    @{ Header = 11; Rows = @(12, 13, 14, 15) },
    @{ Header = 16; Rows = @(17, 18, 19, 20) },
    @{ Header = 21; Rows = @(22, 23) },
    @{ Header = 24; Rows = @(25, 26, 27, 28, 29) },
    @{ Header = 30; Rows = @(31, 32, 33, 34, 35) },
    @{ Header = 36; Rows = @(37, 38, 39) }
)

foreach ($group in $groups) {
    $headerCell = $ws.Cells.Item($group.Header, 1)
    $districtName = $headerCell.Value2

    foreach ($r in $group.Rows) {
        $targetCell = $ws.Cells.Item($r, 1)

        # Copy the header cell's formatting (style/alignment) down onto the
        # target cell, then fill in the matching district text.
        $headerCell.Copy()
        $targetCell.PasteSpecial($xlPasteFormats)
        $targetCell.Value = $districtName
    }
}

$excel.CutCopyMode = $false

# The last group's rows (37-39) use a wrap-text style; once they carry the
# long "Alaska Washington and Oregon..." text the row needs to grow to fit.
$ws.Rows.Item(37).RowHeight = 51
$ws.Rows.Item(38).RowHeight = 51
$ws.Rows.Item(39).RowHeight = 51

# Update the sheet view: scroll so row 32 is at the top and select G38.
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("G38").Select()
